$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update COVID figures (provincias/paises data refresh) ---

# Row 32: Austria - refreshed case counts
$row32 = @(15597, 39, 13228, 1771, 114, 2, 598)
for ($i = 0; $i -lt $row32.Length; $i++) {
    $ws.Cells.Item(32, 2 + $i).Value = $row32[$i]
}

# Row 43: Filipinas - refreshed case counts
$row43 = @(9223, 295, 1214, 7402, 31, 4, 607)
for ($i = 0; $i -lt $row43.Length; $i++) {
    $ws.Cells.Item(43, 2 + $i).Value = $row43[$i]
}

# Row 86 / 87: Lituania overtakes Eslovaquia in total cases, so the two
# countries swap places while the row positions stay fixed (the sheet
# remains sorted descending by "Casos totales"). Row 86 now shows
# Lituania's updated figures, row 87 shows Eslovaquia's.
$ws.Cells.Item(86, 1).Value = "Lituania"
$row86 = @(1410, 4, 635, 729, 17, 0, 46)
for ($i = 0; $i -lt $row86.Length; $i++) {
    $ws.Cells.Item(86, 2 + $i).Value = $row86[$i]
}

$ws.Cells.Item(87, 1).Value = "Eslovaquia"
$row87 = @(1408, 1, 619, 765, 7, 0, 24)
for ($i = 0; $i -lt $row87.Length; $i++) {
    $ws.Cells.Item(87, 2 + $i).Value = $row87[$i]
}
